$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030562134411285
$ws.Range("D2").Value = 1.033206538173445
$ws.Range("E2").Value = 1.038875712721937
$ws.Range("F2").Value = 1.046786514459899
$ws.Range("I2").Value = 1.02359499962809
$ws.Range("J2").Value = 1.035702694217099
$ws.Range("K2").Value = 1.036009554240408
$ws.Range("L2").Value = 1.041662500047411
$ws.Range("M2").Value = 1.049550974837245
$ws.Range("N2").Value = 1.037173510839646
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03215381108919
$ws.Range("D3").Value = 1.03471179359483
$ws.Range("E3").Value = 1.04033598946218
$ws.Range("F3").Value = 1.048434650840371
$ws.Range("I3").Value = 1.023504579208684
$ws.Range("J3").Value = 1.036932603174502
$ws.Range("K3").Value = 1.0373218830727
$ws.Range("L3").Value = 1.042931164785066
$ws.Range("M3").Value = 1.051008653747325
$ws.Range("N3").Value = 1.038405166408842
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033182472109807
$ws.Range("D4").Value = 1.035684889488325
$ws.Range("E4").Value = 1.041279939445195
$ws.Range("F4").Value = 1.049500404049909
$ws.Range("I4").Value = 1.023443335729026
$ws.Range("J4").Value = 1.037726867828675
$ws.Range("K4").Value = 1.038169651329703
$ws.Range("L4").Value = 1.043750624640979
$ws.Range("M4").Value = 1.051950700288228
$ws.Range("N4").Value = 1.039200559009927
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033614629190544
$ws.Range("D5").Value = 1.036093770649095
$ws.Range("E5").Value = 1.041676557336047
$ws.Range("F5").Value = 1.049948288000721
$ws.Range("I5").Value = 1.023416932628374
$ws.Range("J5").Value = 1.038060408364389
$ws.Range("K5").Value = 1.03852572640306
$ws.Range("L5").Value = 1.044094785471043
$ws.Range("M5").Value = 1.052346465061967
$ws.Range("N5").Value = 1.039534573211459
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033687173422782
$ws.Range("D6").Value = 1.036162411617872
$ws.Range("E6").Value = 1.041743138632082
$ws.Range("F6").Value = 1.05002348064651
$ws.Range("I6").Value = 1.023412460900047
$ws.Range("J6").Value = 1.038116389898298
$ws.Range("K6").Value = 1.03858549403571
$ws.Range("L6").Value = 1.044152551856734
$ws.Range("M6").Value = 1.052412900023368
$ws.Range("N6").Value = 1.039590634245568
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033188247747824
$ws.Range("D7").Value = 1.035690353785498
$ws.Range("E7").Value = 1.041285239921015
$ws.Range("F7").Value = 1.049506389307971
$ws.Range("I7").Value = 1.02344298551034
$ws.Range("J7").Value = 1.037731326050824
$ws.Range("K7").Value = 1.038174410495421
$ws.Range("L7").Value = 1.043755224660984
$ws.Range("M7").Value = 1.051955989572085
$ws.Range("N7").Value = 1.039205023563263
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031100315779474
$ws.Range("D8").Value = 1.033715438695404
$ws.Range("E8").Value = 1.039369420912191
$ws.Range("F8").Value = 1.047343660600443
$ws.Range("I8").Value = 1.023565007130014
$ws.Range("J8").Value = 1.036118676871632
$ws.Range("K8").Value = 1.036453355580747
$ws.Range("L8").Value = 1.042091556089259
$ws.Range("M8").Value = 1.050043851354095
$ws.Range("N8").Value = 1.037590084237263
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027411026080918
$ws.Range("D9").Value = 1.030228081309539
$ws.Range("E9").Value = 1.03598589458151
$ws.Range("F9").Value = 1.043526873432085
$ws.Range("I9").Value = 1.023759139049229
$ws.Range("J9").Value = 1.033264644288842
$ws.Range("K9").Value = 1.033409603510073
$ws.Range("L9").Value = 1.039148518373584
$ws.Range("M9").Value = 1.046665091508326
$ws.Range("N9").Value = 1.034731998600844
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024944108354035
$ws.Range("D10").Value = 1.027897741601418
$ws.Range("E10").Value = 1.033724593018178
$ws.Range("F10").Value = 1.040977882752944
$ws.Range("I10").Value = 1.023874603872166
$ws.Range("J10").Value = 1.031353204922555
$ws.Range("K10").Value = 1.031372540575155
$ws.Range("L10").Value = 1.037178341101314
$ws.Range("M10").Value = 1.044405789712117
$ws.Range("N10").Value = 1.032817844771411
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023874016672766
$ws.Range("D11").Value = 1.026887272776325
$ws.Range("E11").Value = 1.032743979980414
$ws.Range("F11").Value = 1.039872953408708
$ws.Range("I11").Value = 1.023921309042974
$ws.Range("J11").Value = 1.030523353244349
$ws.Range("K11").Value = 1.030488490769379
$ws.Range("L11").Value = 1.036323196263314
$ws.Range("M11").Value = 1.04342575841341
$ws.Range("N11").Value = 1.031986814608636
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02347623975599
$ws.Range("D12").Value = 1.0265117167408
$ws.Range("E12").Value = 1.03237950838875
$ws.Range("F12").Value = 1.039462342069674
$ws.Range("I12").Value = 1.023938164268129
$ws.Range("J12").Value = 1.030214772162489
$ws.Range("K12").Value = 1.030159807381616
$ws.Range("L12").Value = 1.036005241231161
$ws.Range("M12").Value = 1.043061459799598
$ws.Range("N12").Value = 1.031677795306246
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023561577894622
$ws.Range("D13").Value = 1.026592285036534
$ws.Range("E13").Value = 1.032457699267964
$ws.Range("F13").Value = 1.039550428436437
$ws.Range("I13").Value = 1.023934571058996
$ws.Range("J13").Value = 1.03028097925438
$ws.Range("K13").Value = 1.030230325167483
$ws.Range("L13").Value = 1.03607345811585
$ws.Range("M13").Value = 1.043139615531749
$ws.Range("N13").Value = 1.031744096419801
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023841142440326
$ws.Range("D14").Value = 1.026856233803733
$ws.Range("E14").Value = 1.032713857357641
$ws.Range("F14").Value = 1.039839016116396
$ws.Range("I14").Value = 1.023922712353276
$ws.Range("J14").Value = 1.030497852766341
$ws.Range("K14").Value = 1.03046132803334
$ws.Range("L14").Value = 1.036296920502645
$ws.Range("M14").Value = 1.04339565097078
$ws.Range("N14").Value = 1.031961277917025
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024013351645638
$ws.Range("D15").Value = 1.027018831506691
$ws.Range("E15").Value = 1.03287165427383
$ws.Range("F15").Value = 1.04001679881495
$ws.Range("I15").Value = 1.023915340507297
$ws.Range("J15").Value = 1.030631430744718
$ws.Range("K15").Value = 1.030603615442074
$ws.Range("L15").Value = 1.036434560891633
$ws.Range("M15").Value = 1.043553366593535
$ws.Range("N15").Value = 1.032095045591453
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02501508578135
$ws.Range("D16").Value = 1.027964772385559
$ws.Range("E16").Value = 1.033789641602625
$ws.Range("F16").Value = 1.041051187050346
$ws.Range("I16").Value = 1.023871434968932
$ws.Range("J16").Value = 1.031408232641062
$ws.Range("K16").Value = 1.031431169375307
$ws.Range("L16").Value = 1.03723505040118
$ws.Range("M16").Value = 1.044470793666129
$ws.Range("N16").Value = 1.032872950635591
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025642930172392
$ws.Range("D17").Value = 1.028557749781807
$ws.Range("E17").Value = 1.034365074542138
$ws.Range("F17").Value = 1.041699702452389
$ws.Range("I17").Value = 1.023843013862359
$ws.Range("J17").Value = 1.031894908533826
$ws.Range("K17").Value = 1.031949733322233
$ws.Range("L17").Value = 1.037736622525416
$ws.Range("M17").Value = 1.045045798800249
$ws.Range("N17").Value = 1.033360317663944
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026008958657824
$ws.Range("D18").Value = 1.028903487495892
$ws.Range("E18").Value = 1.034700575159188
$ws.Range("F18").Value = 1.042077855533128
$ws.Range("I18").Value = 1.023826118601784
$ws.Range("J18").Value = 1.032178568048326
$ws.Range("K18").Value = 1.032252011703005
$ws.Range("L18").Value = 1.038028984345275
$ws.Range("M18").Value = 1.045381022494943
$ws.Range("N18").Value = 1.033644380007492
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026133734307175
$ws.Range("D19").Value = 1.029021352364181
$ws.Range("E19").Value = 1.034814948799851
$ws.Range("F19").Value = 1.042206776794717
$ws.Range("I19").Value = 1.023820303822167
$ws.Range("J19").Value = 1.032275253207596
$ws.Range("K19").Value = 1.032355048667241
$ws.Range("L19").Value = 1.038128639072815
$ws.Range("M19").Value = 1.04549529707875
$ws.Range("N19").Value = 1.033741202470777
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025575587375706
$ws.Range("D20").Value = 1.028494143056313
$ws.Range("E20").Value = 1.03430335055181
$ws.Range("F20").Value = 1.041630134861107
$ws.Range("I20").Value = 1.023846096030143
$ws.Range("J20").Value = 1.031842714632577
$ws.Range("K20").Value = 1.031894116161543
$ws.Range("L20").Value = 1.037682828994939
$ws.Range("M20").Value = 1.044984123545516
$ws.Range("N20").Value = 1.033308049641369
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023758825949247
$ws.Range("D21").Value = 1.026778513706073
$ws.Range("E21").Value = 1.032638431564002
$ws.Range("F21").Value = 1.039754039610723
$ws.Range("I21").Value = 1.02392621804881
$ws.Range("J21").Value = 1.030433998302629
$ws.Range("K21").Value = 1.030393312021792
$ws.Range("L21").Value = 1.036231125199006
$ws.Range("M21").Value = 1.043320262487956
$ws.Range("N21").Value = 1.03189733277265
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022614830413417
$ws.Range("D22").Value = 1.02569853460044
$ws.Range("E22").Value = 1.031590305474058
$ws.Range("F22").Value = 1.038573351184446
$ws.Range("I22").Value = 1.023973741382756
$ws.Range("J22").Value = 1.029546326963721
$ws.Range("K22").Value = 1.029447910689423
$ws.Range("L22").Value = 1.035316547709551
$ws.Range("M22").Value = 1.042272551957323
$ws.Range("N22").Value = 1.031008400838624
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023221451602035
$ws.Range("D23").Value = 1.026271178097183
$ws.Range("E23").Value = 1.032146065846888
$ws.Range("F23").Value = 1.039199365696135
$ws.Range("I23").Value = 1.023948818288664
$ws.Range("J23").Value = 1.030017086655979
$ws.Range("K23").Value = 1.029949258237038
$ws.Range("L23").Value = 1.035801559348375
$ws.Range("M23").Value = 1.042828115832001
$ws.Range("N23").Value = 1.031479829063642
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025606017250857
$ws.Range("D24").Value = 1.028522884617911
$ws.Range("E24").Value = 1.034331241397906
$ws.Range("F24").Value = 1.041661569816517
$ws.Range("I24").Value = 1.023844704313072
$ws.Range("J24").Value = 1.031866299459807
$ws.Range("K24").Value = 1.031919247753341
$ws.Range("L24").Value = 1.037707136583794
$ws.Range("M24").Value = 1.045011992454706
$ws.Range("N24").Value = 1.033331667961758
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028366053722842
$ws.Range("D25").Value = 1.031130566094371
$ws.Range("E25").Value = 1.03686157032531
$ws.Range("F25").Value = 1.044514348941245
$ws.Range("I25").Value = 1.02371141742819
$ws.Range("J25").Value = 1.034003988974582
$ws.Range("K25").Value = 1.034197842754595
$ws.Range("L25").Value = 1.039910767128192
$ws.Range("M25").Value = 1.047539741011017
$ws.Range("N25").Value = 1.035472393240841

Write-Output "Updated 380 kV case values"
